# Se procesan de nuevo los datos con las nuevas dimensiones curadas
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B (municipio) metadata rows
$ws.Range("B2").Value = "sdmx-dimension:refArea"
$ws.Range("B3").Value = "dim"
$ws.Range("B4").Value = "URI-Municipio"

# Update column D (aragon / comunidad) metadata rows
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("D4").Value = "URI-Comunidad"

# Remove the now-obsolete row 5 (mapping-aragon.xlsx)
$ws.Rows("5").Delete()
